$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# ---------------------------------------------------------------------------
# Final lab data (replacing the earlier placeholder numbers) for both result
# tables: "Factor de Carga (PROBING)" (rows 3:5) and
# "Factor de Carga (CHAINING)" (rows 10:12).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 904770.59400000004
$ws.Range("C3").Value = 27672.289000000001
$ws.Range("B4").Value = 800110.14800000004
$ws.Range("C4").Value = 25714.508999999998
$ws.Range("B5").Value = 690105.18799999997
$ws.Range("C5").Value = 24138.222000000002

$ws.Range("B10").Value = 755623.223
$ws.Range("C10").Value = 22031.815999999999
$ws.Range("B11").Value = 737180.34
$ws.Range("C11").Value = 21852.317999999999
$ws.Range("B12").Value = 734420.277
$ws.Range("C12").Value = 23446.78

# ---------------------------------------------------------------------------
# The old "0.00" number formatting is dropped in favour of General, matching
# the big, non-decimal magnitudes of the final measurements.
# ---------------------------------------------------------------------------
$dataRanges = $ws.Range("B3:C5,B10:C12")
$dataRanges.NumberFormat = "General"

# ---------------------------------------------------------------------------
# Light grey shading highlights the first row of each result block ...
# ---------------------------------------------------------------------------
$ws.Range("B3:C3").Interior.Color = 14277081
$ws.Range("B10:C10").Interior.Color = 14277081

# ... the last row of each block keeps the shading too and additionally gets
# a medium black bottom border (closing off the table visually) ...
$ws.Range("B5:C5").Interior.Color = 14277081
$ws.Range("B12:C12").Interior.Color = 14277081

$lastRows = $ws.Range("B5:C5,B12:C12")
$lastRows.Borders.Item(9).LineStyle = 1
$lastRows.Borders.Item(9).Color = 0
$lastRows.Borders.Item(9).Weight = -4138

# ... and the right-hand column ("Consumo de Datos [kB]") of the first and
# last row is additionally set to an explicit black font colour.
$ws.Range("C3").Font.Color = 0
$ws.Range("C10").Font.Color = 0
$ws.Range("B5:C5").Font.Color = 0
$ws.Range("B12:C12").Font.Color = 0

# The middle row of each block (row 4 / row 11) goes back to an unshaded
# cell, only keeping the centred/wrapped alignment.
$ws.Range("B4:C4,B11:C11").Interior.Pattern = -4142

# The last row of each table is taller to accommodate the new bottom border.
$ws.Rows("5").RowHeight = 15.75
$ws.Rows("12").RowHeight = 15.75

# Selection moves to the second table's data, matching the final save state.
$ws.Range("B10:C12").Select()
